$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill row 2 with values 1,2,3,4,5,6 across columns A:F
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 6

# Fill row 3 with values 1,2,3,4,5,6 across columns A:F
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 6

# Update the selected/active cell to G3, matching the post-edit selection
$ws.Range("G3").Select()
